$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.116.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.926.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.43%  '

$ws.Range("E4").Value = '  -0.56%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4730'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4070'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.99'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08445'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.048'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.932.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.520'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.102'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.62%  '

$ws.Range("E16").Value = '  -0.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001068'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.70%  '

$ws.Range("E19").Value = '  -2.16%  '

$ws.Range("E20").Value = '  -6.24%  '

$ws.Range("E21").Value = '  -0.50%  '

$ws.Range("E22").Value = '  -3.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.120.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.52%  '

$ws.Range("E24").Value = '  -4.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.277'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.162.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.32'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -3.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.157'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.716'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9771'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09610'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.447'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.12%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.560'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.69%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.640'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.56%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.059'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02322'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06176'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.240'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6175'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.08%  '

$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1908'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5904'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.286'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.038'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.07%  '

$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06812'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.97%  '
